# Limits-V5.xlsx update — "22 data calculated upto boring no 87 and D-20"
#
# The "main" sheet pulls its borehole hydrometer-analysis inputs from an
# external workbook (Hydrometer_V8.xlsm, sheet "input-output") via cached
# formulas such as ='[1]input-output'!$E$24. That source file isn't part of
# this workbook and isn't reachable from here, so — exactly like a user who
# re-keys the refreshed numbers by hand after the linked file changed — we
# push the new readings straight onto the cells that used to hold those
# external-reference formulas. Every other changed cell in the workbook
# (rows 20-22, 39-41, 47-49, 54, 56 on "main"; the mirrored cells and chart
# source cells on "Report") is a local formula that recalculates on its own
# once these inputs change.

$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("main")

# --- Row 16: Sieve/finer-than data (was '[1]input-output'!$D$24:$H$24) ---
# C16 (=D24) stays 15 — unchanged in the source commit.
$main.Range("D16").Value2 = 19
$main.Range("E16").Value2 = 24
$main.Range("F16").Value2 = 29
$main.Range("G16").Value2 = 34

# --- Row 17 (was '[1]input-output'!$D$25:$H$25) ---
$main.Range("C17").Value2 = 10.16
$main.Range("D17").Value2 = 8.96
$main.Range("E17").Value2 = 9.83
$main.Range("F17").Value2 = 9.7
$main.Range("G17").Value2 = 10.23

# --- Row 18 (was '[1]input-output'!$D$26:$H$26) ---
$main.Range("C18").Value2 = 44.56
$main.Range("D18").Value2 = 37.31
$main.Range("E18").Value2 = 44.85
$main.Range("F18").Value2 = 38.55
$main.Range("G18").Value2 = 40.84

# --- Row 19 (was '[1]input-output'!$D$27:$H$27) ---
$main.Range("C19").Value2 = 35.09
$main.Range("D19").Value2 = 29.44
$main.Range("E19").Value2 = 35.48
$main.Range("F19").Value2 = 30.44
$main.Range("G19").Value2 = 32.79

# --- Row 44 (was '[1]input-output'!$M$24:$O$24) ---
$main.Range("D44").Value2 = 7.02
$main.Range("E44").Value2 = 11.18
$main.Range("F44").Value2 = 9.09

# --- Row 45 (was '[1]input-output'!$M$25:$O$25) ---
$main.Range("D45").Value2 = 45.45
$main.Range("E45").Value2 = 47.74
$main.Range("F45").Value2 = 44.72

# --- Row 46 (was '[1]input-output'!$M$26:$O$26) ---
$main.Range("D46").Value2 = 38
$main.Range("E46").Value2 = 40.66
$main.Range("F46").Value2 = 37.86

# Everything downstream (C20:G22 grain-size differences/percentages, the
# TREND-based D10/Cu cells G39:G41, the Atterberg rows 47-49, G54/G56 on
# "main", and all of their mirrors + G-column copies on "Report") is a plain
# in-workbook formula, so it re-evaluates automatically on recalc — no
# further writes are needed for those.
$excel.CalculateFull()

# --- View state ---
# "main": scroll position moves to row 40 while the selection stays on G20.
$main.Activate()
$main.Range("G20").Select()
$mainWin = $excel.ActiveWindow
$mainWin.ScrollRow = 40
$mainWin.ScrollColumn = 1

# "Report": scroll position moves to row 16 and the selection moves to M25.
$report = $wb.Worksheets.Item("Report")
$report.Activate()
$report.Range("M25").Select()
$reportWin = $excel.ActiveWindow
$reportWin.ScrollRow = 16
$reportWin.ScrollColumn = 3
